$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (some with trailing zeros, e.g. "32.50") are preserved exactly as text,
# matching the workbook's existing inline-string convention.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '41.778.59'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '2.265.53'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '303.69'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").Value = '92.21'
$ws.Range("E6").Value = '  +1.51%  '
$ws.Range("E7").Value = '  +1.98%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '32.50'
$ws.Range("E10").Value = '  +2.19%  '
$ws.Range("D11").Value = '53.49'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '2.612.46'
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '2.266.42'
$ws.Range("E17").Value = '  -5.77%  '
$ws.Range("D18").Value = '0.770'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").Value = '41.670.37'
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("D20").Value = '12.43'
$ws.Range("E20").Value = '  +4.53%  '
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("E22").Value = '  +1.63%  '
$ws.Range("D23").Value = '67.18'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").Value = '239.80'
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("E25").Value = '  +1.31%  '
$ws.Range("D27").Value = '1.93'
$ws.Range("E27").Value = '  +4.11%  '
$ws.Range("D28").Value = '23.95'
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = '35.42'
$ws.Range("E31").Value = '  +6.64%  '
$ws.Range("D32").Value = '160.42'
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").Value = '5.25'
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = '0.0743'
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  +3.13%  '
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").Value = '2.016.49'
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").Value = '19.30'
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("D46").Value = '10.33'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("E47").Value = '  +6.36%  '
$ws.Range("D48").Value = '2.89'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").Value = '1.16'
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("D51").Value = '52.32'
$ws.Range("E51").Value = '  +3.24%  '
